$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F4").Value = 1.58
$ws.Range("G4").Value = 1.65
$ws.Range("I4").Value = 7.2
$ws.Range("P4").Value = 2.02
$ws.Range("Q4").Value = 1.87
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 1.8
$ws.Range("K5").Value = 3.95
$ws.Range("N5").Value = 3.05
$ws.Range("P5").Value = 1.73
$ws.Range("Q5").Value = 2.14
$ws.Range("R5").Value = 1.26
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 1.87
$ws.Range("U5").Value = 1.84
$ws.Range("F6").Value = 1.42
$ws.Range("G6").Value = 1.53
$ws.Range("H6").Value = 5.9
$ws.Range("I6").Value = 13
$ws.Range("K6").Value = 5.4
$ws.Range("P6").Value = 1.74
$ws.Range("F8").Value = 1.44
$ws.Range("P9").Value = 1.88
$ws.Range("Q9").Value = 1.74
$ws.Range("P10").Value = 1.63
$ws.Range("U10").Value = 1.92
$ws.Range("AJ10").Value = 50
$ws.Range("F11").Value = 3.95
$ws.Range("G11").Value = 6.6
$ws.Range("H11").Value = 1.79
$ws.Range("I11").Value = 2.04
$ws.Range("J11").Value = 3.2
$ws.Range("K11").Value = 3.95
$ws.Range("P11").Value = 1.59
$ws.Range("Q11").Value = 2.02
$ws.Range("G15").Value = 2.62
$ws.Range("H15").Value = 2.66
$ws.Range("J15").Value = 2.84
$ws.Range("K15").Value = 4.1
$ws.Range("F16").Value = 1.77
$ws.Range("H16").Value = 3.2
$ws.Range("J16").Value = 2.9
$ws.Range("K16").Value = 8.800000000000001
$ws.Range("F17").Value = 3.6
$ws.Range("G17").Value = 4.6
$ws.Range("H17").Value = 1.87
$ws.Range("I17").Value = 2.06
$ws.Range("J17").Value = 3.95
$ws.Range("K17").Value = 5.1
$ws.Range("P17").Value = 2.5
$ws.Range("Q17").Value = 1.54
$ws.Range("F19").Value = 2.48
$ws.Range("J19").Value = 2.98
$ws.Range("P19").Value = 1.57
$ws.Range("Q19").Value = 2.54
$ws.Range("G20").Value = 3.2
$ws.Range("H20").Value = 2.76
$ws.Range("K20").Value = 3.65
$ws.Range("P20").Value = 1.67
$ws.Range("T21").Value = 1.99
$ws.Range("Z21").Value = 19
$ws.Range("AB21").Value = 9.4
$ws.Range("AL22").Value = 25
$ws.Range("Y23").Value = 11.5
$ws.Range("AB23").Value = 8.800000000000001
$ws.Range("AE23").Value = 46
$ws.Range("F24").Value = 2.76
$ws.Range("G24").Value = 2.98
$ws.Range("H24").Value = 2.72
$ws.Range("I24").Value = 2.92
$ws.Range("K24").Value = 3.55
$ws.Range("L24").Value = 1.4
$ws.Range("P24").Value = 1.93
$ws.Range("Q24").Value = 1.94
$ws.Range("S24").Value = 3.35
$ws.Range("T24").Value = 1.72
$ws.Range("V24").Value = 1.52
$ws.Range("W24").Value = 1.51
$ws.Range("X24").Value = 17
$ws.Range("Y24").Value = 13.5
$ws.Range("G25").Value = 2.5
$ws.Range("H25").Value = 3.55
$ws.Range("I25").Value = 4.3
$ws.Range("J25").Value = 3.1
$ws.Range("K25").Value = 3.6
$ws.Range("L25").Value = 1.48
$ws.Range("M25").Value = 1.09
$ws.Range("N25").Value = 2.9
$ws.Range("P25").Value = 1.65
$ws.Range("Q25").Value = 2.22
$ws.Range("S25").Value = 3.85
$ws.Range("V25").Value = 1.31
$ws.Range("W25").Value = 1.69
$ws.Range("Y25").Value = 14.5
$ws.Range("AA25").Value = 100
$ws.Range("AB25").Value = 9.800000000000001
$ws.Range("AF25").Value = 16.5
$ws.Range("AO25").Value = 85
